$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.895.27'
$ws.Range("E2").Value = '  -0.29%  '

# Row 3
$ws.Range("D3").Value = '2.213.45'
$ws.Range("E3").Value = '  -1.66%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '''255.11'
$ws.Range("E5").Value = '  +3.99%  '

# Row 6
$ws.Range("D6").Value = '''0.616'

# Row 7
$ws.Range("D7").Value = '''76.22'
$ws.Range("E7").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").Value = '''0.595'
$ws.Range("E9").Value = '  -3.58%  '

# Row 10
$ws.Range("D10").Value = '''41.76'
$ws.Range("E10").Value = '  +1.80%  '

# Row 11
$ws.Range("D11").Value = '''0.0908'
$ws.Range("E11").Value = '  -3.03%  '

# Row 12
$ws.Range("D12").Value = '''6.92'
$ws.Range("E12").Value = '  -0.89%  '

# Row 14
$ws.Range("D14").Value = '2.544.82'
$ws.Range("E14").Value = '  -1.65%  '

# Row 15
$ws.Range("D15").Value = '''14.38'
$ws.Range("E15").Value = '  -1.76%  '

# Row 16
$ws.Range("D16").Value = '2.214.78'
$ws.Range("E16").Value = '  -1.80%  '

# Row 17
$ws.Range("D17").Value = '''0.781'
$ws.Range("E17").Value = '  -3.39%  '

# Row 18
$ws.Range("D18").Value = '42.837.91'
$ws.Range("E18").Value = '  -0.24%  '

# Row 19
$ws.Range("E19").Value = '  -2.98%  '

# Row 20
$ws.Range("D20").Value = '''71.29'
$ws.Range("E20").Value = '  +0.12%  '

# Row 21
$ws.Range("D21").Value = '''5.94'
$ws.Range("E21").Value = '  -0.93%  '

# Row 22
$ws.Range("D22").Value = '''229.96'
$ws.Range("E22").Value = '  -0.46%  '

# Row 23
$ws.Range("D23").Value = '''2.20'
$ws.Range("E23").Value = '  +0.07%  '

# Row 24
$ws.Range("D24").Value = '''9.19'
$ws.Range("E24").Value = '  -8.47%  '

# Row 25
$ws.Range("E25").Value = '  -0.11%  '

# Row 26
$ws.Range("D26").Value = '''10.66'
$ws.Range("E26").Value = '  -2.05%  '

# Row 27
$ws.Range("D27").Value = '''40.55'
$ws.Range("E27").Value = '  +3.89%  '

# Row 28
$ws.Range("E28").Value = '  -2.80%  '

# Row 29
$ws.Range("E29").Value = '  +0.03%  '

# Row 30
$ws.Range("E30").Value = '  -3.35%  '

# Row 31
$ws.Range("D31").Value = '''173.73'
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("D32").Value = '''20.23'
$ws.Range("E32").Value = '  -0.57%  '

# Row 33
$ws.Range("D33").Value = '''0.0845'
$ws.Range("E33").Value = '  +5.91%  '

# Row 34
$ws.Range("D34").Value = '''5.21'
$ws.Range("E34").Value = '  -2.43%  '

# Row 35
$ws.Range("E35").Value = '  -1.18%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '''0.0353'
$ws.Range("E36").Value = '  +6.00%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.107'
$ws.Range("E37").Value = '  -2.89%  '

# Row 38
$ws.Range("D38").Value = '''4.30'
$ws.Range("E38").Value = '  -1.23%  '

# Row 39
$ws.Range("D39").Value = '''12.56'
$ws.Range("E39").Value = '  -3.26%  '

# Row 40
$ws.Range("E40").Value = '  -2.37%  '

# Row 41
$ws.Range("D41").Value = '''2.74'
$ws.Range("E41").Value = '  +14.35%  '

# Row 42
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").Value = '''5.27'
$ws.Range("E42").Value = '  -5.14%  '

# Row 43
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.198'
$ws.Range("E43").Value = '  -2.60%  '

# Row 44
$ws.Range("D44").Value = '''59.70'
$ws.Range("E44").Value = '  -0.84%  '

# Row 45
$ws.Range("D45").Value = '''101.84'
$ws.Range("E45").Value = '  -3.78%  '

# Row 46
$ws.Range("D46").Value = '''0.0979'
$ws.Range("E46").Value = '  -1.59%  '

# Row 47
$ws.Range("D47").Value = '''8.27'
$ws.Range("E47").Value = '  -4.97%  '

# Row 48
$ws.Range("D48").Value = '''0.453'
$ws.Range("E48").Value = '  -1.38%  '

# Row 49
$ws.Range("D49").Value = '''1.11'
$ws.Range("E49").Value = '  -0.50%  '

# Row 50
$ws.Range("E50").Value = '  -1.34%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.436.51'
$ws.Range("E51").Value = '  -0.94%  '
